$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text (e.g. "48.168.73",
# "0.555") rather than numbers. A handful of the new prices parse as valid
# numbers (e.g. "109.91"), so Excel would silently coerce them to the Number
# type on assignment. Pre-format just those cells as Text ("@") first so the
# written value is preserved verbatim as a string, matching the rest of the column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '48.169.12'
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('D3').Value = '2.519.62'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '109.91'
$ws.Range('E5').Value = '  +1.10%  '
$ws.Range('D6').Value = '322.73'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('E7').Value = '  +2.05%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.554'
$ws.Range('E9').Value = '  +4.05%  '
$ws.Range('D10').Value = '40.62'
$ws.Range('E10').Value = '  +4.83%  '
$ws.Range('D11').Value = '20.53'
$ws.Range('E11').Value = '  +12.55%  '
$ws.Range('D12').Value = '0.0826'
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('E13').Value = '  +1.27%  '
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('D15').Value = '2.915.30'
$ws.Range('E15').Value = '  +1.30%  '
$ws.Range('D16').Value = '2.525.76'
$ws.Range('D17').Value = '0.854'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').Value = '48.004.72'
$ws.Range('E18').Value = '  +2.02%  '
$ws.Range('D19').Value = '13.24'
$ws.Range('E19').Value = '  +4.23%  '
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').Value = '0.0₃0948'
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('E22').Value = '  -1.99%  '
$ws.Range('D23').Value = '72.11'
$ws.Range('E23').Value = '  +2.02%  '
$ws.Range('D24').Value = '264.93'
$ws.Range('E24').Value = '  +7.80%  '
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '26.12'
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('D28').Value = '10.15'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('E29').Value = '  +3.58%  '
$ws.Range('D30').Value = '36.40'
$ws.Range('E30').Value = '  +3.80%  '
$ws.Range('D31').Value = '2.21'
$ws.Range('E31').Value = '  -2.91%  '
$ws.Range('D32').Value = '49.74'
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('D33').Value = '19.89'
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '0.0791'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('E37').Value = '  +1.49%  '
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('E39').Value = '  +1.68%  '
$ws.Range('D40').Value = '0.112'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '22.13'
$ws.Range('E41').Value = '  +3.06%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = '120.13'
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '2.20'
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range('E44').Value = '  +2.21%  '
$ws.Range('D45').Value = '2.019.75'
$ws.Range('E45').Value = '  +2.06%  '
$ws.Range('D46').Value = '3.17'
$ws.Range('E46').Value = '  +5.05%  '
$ws.Range('E47').Value = '  +8.72%  '
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('D50').Value = '5.25'
$ws.Range('E50').Value = '  +2.54%  '
$ws.Range('D51').Value = '78.87'
$ws.Range('E51').Value = '  +2.65%  '
